$d = $word.ActiveDocument

# 1) Straightforward whole-paragraph text replacements (date, title, body paragraphs, URL)
$null = $d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק 29.06.24:⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק 28.06.24:⚡️🚀", 2)
$null = $d.Content.Find.Execute("What Are the Odds? Language Models Are Capable of Probabilistic Reasoning", $true, $false, $false, $false, $false, $true, 1, $false, "On-Policy Distillation OF LANGUAGE MODELS: LEARNING FROM SELF-GENERATED MISTAKES", 2)
$null = $d.Content.Find.Execute("הסקירה הזו הולכת להיות ממש קצרה. לפני ימיים (27.06) סקרתי מאמר שבדק האם מודלי שפה ענקיים מסוגלים לבצע רגרסיה לוגיסטית והגיע למסקנה שבלי עזרה ורמזים מאוד משמעתיים הם לא מצליחים לפתור אותה. ", $true, $false, $false, $false, $false, $true, 1, $false, "מזמן לא סקרתי מאמר על שיטות זיקוק של ידע(knowledge distillation) - לא נתקלתי במאמרים מגניבים בנושא המעניין הזה. מה זה זיקוק ידע ממודל גדול למודל קטן יותר? למעשה זה ניסיון להעתיק למודל הקטן את הידע שיש למודל הגדול כלומר לגרום לו להפגין ביצועים הדומים למודל הגדול.", 2)
$null = $d.Content.Find.Execute("הפעם המחברים בדקו האם מודלי שפה מסוגלים ״לנתח התפלגויות הסתברותיות״. למשל אומרים למודל שפה שאיזשהו ערך מפולג גאוסית עם תוחלת 3 ושונות 2 ושואלים אותו מה האחוזון ה-80 של ההתפלגות או לאיזה אחוזון שייכת דגימה בעלת ערך 4. באופן די מפתיע המודל מצליח לא רע בשאלות האלו למרות שקיבל הוראה לא להריץ קוד (זה יכול לעזור כמו שאתם מבינים).", $true, $false, $false, $false, $false, $true, 1, $false, "יש כמה שיטות לעשות זאת - הפשוטה ביותר זה לאמן אותו על הדאטה שהמודל הגדול אומן עליה. יש שיטות המאמנות את המודל הקטן על הדאטה המיוצר על ידי המודל הגדול. אם יש לנו גישה להתפלגויות (של הטוקנים) אז מאמנים את המודל הקטן לחקות את התפלגות הטוקנים שהמודל הגדול מוציא. אם יש לנו אקטיבציות של השכבות של המודל הגדול ניתן לנסות לחקות גם אותם (אם המודל הקטן הוא בעל אותה ארכיטקטורה אבל עם פחות שכבות).", 2)
$null = $d.Content.Find.Execute("אז מה לדעתכם קורה כאן? איך המודל מצליח לפתור את השאלות האלו?", $true, $false, $false, $false, $false, $true, 1, $false, "בכל גישות האלו אנו מאמנים (או פיינטיון) את המודל הקטן בצורה supervised רגילה. כלומר יש לנו סט של דוגמאות (ground-truth או שנוצרו על ידי המודל הגדול) אנו מאמנים את המודל הקטן עליהם. המאמר שנסקור היום מציעה להשתמש בגישה מעולמות למידה באמצעות חיזוקים (reinforcement learning) ממשפחת on-policy. זה אומר שהאימון מתבצע על הדוגמאות שהרשת המאומנת עצמה יוצרת במהלך האימון (והיא משתנה כמובן).", 2)
$null = $d.Content.Find.Execute("https://arxiv.org/abs/2406.12830", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2306.13649", 2)

# 2) Insert the new paragraphs that were added before the final URL paragraph.
#    The URL paragraph is always the last paragraph in the document; insert new
#    empty paragraphs directly before it (preserving its own 'Normal' style/formatting)
#    and then fill each with its text, from first to last so ordering is preserved.
$urlPara = $d.Paragraphs.Last
$urlPara.Range.InsertParagraphBefore()
$newIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newIndex).Range.Text = "המאמר הלך צעד אחד קדימה והחליט לשלב את שיטת אימון on-policy יחד עם האימון הסטנדרטי של זיקוק ידע. כלומר בהסתברות alpha השיטה בוחרת דוגמא מדאטהסט האימון ובשאר המקרים היא מגרילה דאטה מהמודל הקטן. כל פעם המודל מנסה למזער את המרחק בין התפלגות הטוקנים של הדוגמא (מהדאטהסט או מהמודל הקטן). "

$urlPara = $d.Paragraphs.Last
$urlPara.Range.InsertParagraphBefore()
$newIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newIndex).Range.Text = "בד״כ כלל המרחק בין התפלגויות של הטוקנים בשיטות זיקוק ידע נמדד על KL divergence סטנדרטי (כלומר forward). המאמר מציע לשכלל את הגישה הזו עקב חולשה שיש ל- forward KL. החולשה הזו קשורה לעובדה ש-forward KL מנסה לקרב את התפלגות המודל המאומן לאזור המוד(mode) של התפלגות היעד (התפלגות המודל הגדול במקרה שלנו. הכוונה כאן שהתפלגות המודל המאומן עלולה ״להתרכז באזור בעל מסה הסתברותיות גבוהה״, נגיד ליד איזה מוד של ההתפלגות ומתעלמת מאיזורים אחרים שיש בהם מסה הסתברותית ליד מודים חלשים יותר של ההתפלגות."

$urlPara = $d.Paragraphs.Last
$urlPara.Range.InsertParagraphBefore()
$newIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newIndex).Range.Text = "למזלנו יש לנו reverse KL שהופך את המונה ואת המכנה בלוג של forward KL. ניתן להראות כי forward KL מנסה ״לכסות״ את כל האזור בה התפלגות היעד גדולה מאפס ובכך משלימה את forward KL. ניתן לשלב אותם לינארית (באופן קמור עם מקדם beta ו- 1-beta) ואז מקבל Jensen Shannon Convergence או JSD שנותן מענה לבעיה האינהרנטית של forward KL. ובה המאמר משתמש במקום forward KL הרגיל. "

$urlPara = $d.Paragraphs.Last
$urlPara.Range.InsertParagraphBefore()
$newIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newIndex).Range.Text = "ניתן לשלב את פונקציית הלוס של המאמר עם עוד איבר האחראי על מקסום פונקציית reward כלשהי עבור המודל הקטן (כמו ב-RLHF)."

$urlPara = $d.Paragraphs.Last
$urlPara.Range.InsertParagraphBefore()
$newIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newIndex).Range.Text = "ושכחתי להגיד(לא קשור למאמר) ש- forward KL זה בדיוק מה יש לנו בכל פונקציית לוס המבוססת על cross entropy (נגיד במשימות סיווג)."

